$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain-looking number (e.g. "1.00", "590.11") must be
# forced to Text format first, otherwise Excel's COM Value setter auto-converts them
# to numeric cells (stripping the literal formatting the source data relies on).
$textForceCells = @(
    "D5",
    "D6",
    "D13",
    "D14",
    "D19",
    "D20",
    "D21",
    "D23",
    "D24",
    "D26",
    "D27",
    "D29",
    "D31",
    "D32",
    "D35",
    "D36",
    "D38",
    "D39",
    "D43",
    "D46",
    "D49",
)
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated Price (D) / Volume(1h) (E) values scraped for this run.
$ws.Range("D2").Value = '66.902.99'
$ws.Range("E2").Value = '  -2.23%  '
$ws.Range("D3").Value = '2.621.54'
$ws.Range("E3").Value = '  -3.45%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = '590.11'
$ws.Range("E5").Value = '  -2.89%  '
$ws.Range("D6").Value = '164.54'
$ws.Range("E6").Value = '  -1.43%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  -3.81%  '
$ws.Range("D9").Value = '2.620.36'
$ws.Range("E9").Value = '  -3.48%  '
$ws.Range("E10").Value = '  -2.95%  '
$ws.Range("E11").Value = '  +1.15%  '
$ws.Range("E12").Value = '  -1.51%  '
$ws.Range("D13").Value = '5.22'
$ws.Range("E13").Value = '  -1.65%  '
$ws.Range("D14").Value = '27.38'
$ws.Range("E14").Value = '  -4.02%  '
$ws.Range("D15").Value = '3.103.10'
$ws.Range("E15").Value = '  -3.27%  '
$ws.Range("E16").Value = '  -3.57%  '
$ws.Range("D17").Value = '66.948.03'
$ws.Range("E17").Value = '  -2.11%  '
$ws.Range("D18").Value = '2.616.35'
$ws.Range("E18").Value = '  -3.35%  '
$ws.Range("D19").Value = '12.06'
$ws.Range("E19").Value = '  +1.53%  '
$ws.Range("D20").Value = '8.01'
$ws.Range("E20").Value = '  +5.07%  '
$ws.Range("D21").Value = '358.20'
$ws.Range("E21").Value = '  -3.48%  '
$ws.Range("E22").Value = '  -3.87%  '
$ws.Range("D23").Value = '4.63'
$ws.Range("E23").Value = '  -6.38%  '
$ws.Range("D24").Value = '10.86'
$ws.Range("E24").Value = '  +7.35%  '
$ws.Range("E25").Value = '  -7.17%  '
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("D27").Value = '71.05'
$ws.Range("E27").Value = '  -2.86%  '
$ws.Range("E28").Value = '  -3.64%  '
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  -0.01%  '
$ws.Range("E30").Value = '  -3.54%  '
$ws.Range("D31").Value = '550.58'
$ws.Range("E31").Value = '  -5.11%  '
$ws.Range("D32").Value = '7.90'
$ws.Range("E32").Value = '  -3.09%  '
$ws.Range("E33").Value = '  -4.47%  '
$ws.Range("E34").Value = '  -3.84%  '
$ws.Range("D35").Value = '0.133'
$ws.Range("E35").Value = '  +1.45%  '
$ws.Range("D36").Value = '0.999'
$ws.Range("E36").Value = '  +0.01%  '
$ws.Range("E37").Value = '  -5.84%  '
$ws.Range("D38").Value = '157.41'
$ws.Range("E38").Value = '  -2.55%  '
$ws.Range("D39").Value = '19.13'
$ws.Range("E39").Value = '  -3.82%  '
$ws.Range("E40").Value = '  -3.24%  '
$ws.Range("E41").Value = '  -3.88%  '
$ws.Range("E42").Value = '  -4.97%  '
$ws.Range("D43").Value = '17.90'
$ws.Range("E43").Value = '  -0.47%  '
$ws.Range("E44").Value = '  -0.03%  '
$ws.Range("E45").Value = '  -6.42%  '
$ws.Range("D46").Value = '40.26'
$ws.Range("E46").Value = '  -1.10%  '
$ws.Range("D47").Value = '0.0₆0297'
$ws.Range("E47").Value = '  -4.81%  '
$ws.Range("E48").Value = '  -2.33%  '
$ws.Range("D49").Value = '151.50'
$ws.Range("E49").Value = '  -2.08%  '
$ws.Range("E50").Value = '  -2.70%  '
$ws.Range("E51").Value = '  -3.37%  '
